{"js": "// Office.js (Word JavaScript API) implementation of:\n//   \"Added some of Nathan's contributions\"\n//\n// Two content edits are made to the document body:\n//   1. In the \"Included files\" section, the README entry is updated from\n//      \"README - The file you are currently reading. ...\" to\n//      \"README.pdf - The file you are currently reading. ...\", reflecting\n//      that the README is actually a PDF file.\n//   2. In the \"Group contributions\" section, Nathan Henninger's paragraph\n//      (which previously only contained a tab character) gets his actual\n//      contribution sentence appended: \"Cleaned up code, fixed bugs, added\n//      input validation, and helped with documentation.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst OLD_README = \"README - The file you are currently reading\";\nconst NEW_README = \"README.pdf - The file you are currently reading\";\nconst NATHAN_LABEL = \"Nathan Henninger:\";\nconst NATHAN_SENTENCE =\n  \"Cleaned up code, fixed bugs, added input validation, and helped with documentation.\";\n\n// --- Edit 1: README -> README.pdf ---------------------------------------\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(OLD_README) !== -1 && p.text.indexOf(NEW_README) === -1) {\n    const updated = p.text.replace(OLD_README, NEW_README);\n    p.getRange().insertText(updated, Word.InsertLocation.replace);\n    break;\n  }\n}\n\n// --- Edit 2: add Nathan Henninger's contribution sentence ---------------\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const isBlankOrTabOnly = p.text.trim() === \"\";\n  const precededByNathanLabel =\n    i > 0 && paragraphs.items[i - 1].text.trim() === NATHAN_LABEL;\n  const alreadyDone = p.text.indexOf(NATHAN_SENTENCE) !== -1;\n\n  if (isBlankOrTabOnly && precededByNathanLabel && !alreadyDone) {\n    p.insertText(NATHAN_SENTENCE, Word.InsertLocation.end);\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop implementation of:\n#   \"Added some of Nathan's contributions\"\n#\n# Two content edits are made to the document:\n#   1. In the \"Included files\" section, the README entry is updated from\n#      \"README - The file you are currently reading. ...\" to\n#      \"README.pdf - The file you are currently reading. ...\", reflecting\n#      that the README is actually a PDF file.\n#   2. In the \"Group contributions\" section, Nathan Henninger's paragraph\n#      (which previously only contained a tab character) gets his actual\n#      contribution sentence appended: \"Cleaned up code, fixed bugs, added\n#      input validation, and helped with documentation.\"\n\n$d = $word.ActiveDocument\n\n$nathanSentence = \"Cleaned up code, fixed bugs, added input validation, and helped with documentation.\"\n\n# --- Edit 1: README -> README.pdf ---------------------------------------\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"README - The file you are currently reading\"\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Replacement.Text = \"README.pdf - The file you are currently reading\"\n$rng.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# --- Edit 2: add Nathan Henninger's contribution sentence ---------------\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.Trim()\n    $prevText = \"\"\n    if ($i -gt 1) {\n        $prevText = $d.Paragraphs.Item($i - 1).Range.Text.Trim()\n    }\n\n    if ($text -eq \"\" -and $prevText -eq \"Nathan Henninger:\" -and $text -notlike \"*$nathanSentence*\") {\n        $p.Range.InsertAfter($nathanSentence)\n        break\n    }\n}\n"}
